# Automatische test-sync: 2025-08-04 20:36:50
$wb = $excel.ActiveWorkbook

$wsLogs = $wb.Worksheets.Item("Logs")
$wsDashboard = $wb.Worksheets.Item("Dashboard")

# Append new row 17 to the "Logs" sheet
$wsLogs.Range("A17").Value = "Wil je deze klant bellen?"
$wsLogs.Range("B17").Value = "mailmind.test@zohomail.eu"
$wsLogs.Range("C17").Value = "Testmail #5: Wil je deze klant bellen?"
$wsLogs.Range("D17").Value = "Klantenservice / Contact"
$wsLogs.Range("E17").Value = "Geachte heer/mevrouw,`nBedankt voor uw e-mail. We zullen de klant zo snel mogelijk contacteren. Mocht u nog meer informatie hebben die u met ons wilt delen, dan horen we dat graag.`nMet vriendelijke groet,`n[Naam bedrijf]"
$wsLogs.Range("F17").Value = "2025-08-04 20:36:10"
$wsLogs.Range("G17").Value = "Ja"
$wsLogs.Range("H17").Value = "Nee"
$wsLogs.Range("I17").Value = "Ja"
$wsLogs.Range("J17").Value = "Nee"

# Undo the auto row-height bump triggered by the multi-line text in E17,
# matching the source row (no explicit height override)
$wsLogs.Rows.Item(17).AutoFit()

# Extend the conditional-formatting ranges so they cover the new row too
$columns = @("D", "G", "H", "I", "J")
foreach ($col in $columns) {
    $oldRange = $wsLogs.Range($col + "2:" + $col + "16")
    $newRange = $wsLogs.Range($col + "2:" + $col + "17")
    $fcs = $oldRange.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRange)
    }
}

# Update the summary count on the "Dashboard" sheet for "Klantenservice / Contact"
$wsDashboard.Range("B6").Value = 2
